$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 34.58258212375859
$ws.Range("I2").Value = 34.91871657754011
$ws.Range("J2").Value = 34.90099312452253
$ws.Range("K2").Value = 34.52941176470589
$ws.Range("L2").Value = 36.01176470588235
$ws.Range("M2").Value = 37.4
$ws.Range("N2").Value = 38.69411764705882
$ws.Range("O2").Value = 39.89411764705883

# Row 4
$ws.Range("H4").Value = 18.06493506493507
$ws.Range("I4").Value = 17.63636363636364
$ws.Range("J4").Value = 16.94805194805195
$ws.Range("K4").Value = 16
$ws.Range("L4").Value = 17
$ws.Range("M4").Value = 18
$ws.Range("N4").Value = 19
$ws.Range("O4").Value = 20

# Row 5
$ws.Range("H5").Value = 34.58258212375859
$ws.Range("I5").Value = 34.91871657754011
$ws.Range("J5").Value = 34.90099312452253
$ws.Range("K5").Value = 34.52941176470589
$ws.Range("L5").Value = 36.01176470588235
$ws.Range("M5").Value = 37.4
$ws.Range("N5").Value = 38.69411764705882
$ws.Range("O5").Value = 39.89411764705883

# Row 7
$ws.Range("H7").Value = 18.06493506493507
$ws.Range("I7").Value = 17.63636363636364
$ws.Range("J7").Value = 16.94805194805195
$ws.Range("K7").Value = 16
$ws.Range("L7").Value = 17
$ws.Range("M7").Value = 18
$ws.Range("N7").Value = 19
$ws.Range("O7").Value = 20
